# The author removed one post entry from the sheet: the row containing
# the "「このラマは誇らしげにして立っている」" (llama) entry, which was
# row 698. Deleting that row shifts every following row up by one,
# which matches the rest of the diff (all subsequent row numbers
# decrease by 1) and the updated dimension (A1:C838 -> A1:C837).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(698).Delete()
